$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'52.027.79"
$ws.Range("E2").Value = "'  +0.72%  "
$ws.Range("D3").Value = "'3.008.89"
$ws.Range("E3").Value = "'  +3.22%  "
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("D5").Value = "'354.33"
$ws.Range("E5").Value = "'  +0.39%  "
$ws.Range("D6").Value = "'106.85"
$ws.Range("E6").Value = "'  -2.95%  "
$ws.Range("E7").Value = "'  -1.83%  "
$ws.Range("E8").Value = "'  +0.01%  "
$ws.Range("D9").Value = "'0.612"
$ws.Range("E9").Value = "'  -2.26%  "
$ws.Range("D10").Value = "'38.08"
$ws.Range("E10").Value = "'  -2.41%  "
$ws.Range("E11").Value = "'  +2.95%  "
$ws.Range("D12").Value = "'0.0856"
$ws.Range("E12").Value = "'  -2.91%  "
$ws.Range("D13").Value = "'18.98"
$ws.Range("E13").Value = "'  -3.13%  "
$ws.Range("D14").Value = "'3.478.94"
$ws.Range("E14").Value = "'  +3.30%  "
$ws.Range("D15").Value = "'7.61"
$ws.Range("E15").Value = "'  -2.93%  "
$ws.Range("D16").Value = "'2.989.45"
$ws.Range("E16").Value = "'  +2.91%  "
$ws.Range("E17").Value = "'  +4.15%  "
$ws.Range("D18").Value = "'52.025.42"
$ws.Range("E18").Value = "'  +0.84%  "
$ws.Range("E19").Value = "'  +4.89%  "
$ws.Range("D20").Value = "'7.49"
$ws.Range("E20").Value = "'  -0.24%  "
$ws.Range("D21").Value = "'13.60"
$ws.Range("E21").Value = "'  -1.70%  "
$ws.Range("D22").Value = "'0.0₃0972"
$ws.Range("E22").Value = "'  -0.50%  "
$ws.Range("D23").Value = "'69.06"
$ws.Range("E23").Value = "'  -2.13%  "
$ws.Range("D24").Value = "'263.68"
$ws.Range("E24").Value = "'  -2.01%  "
$ws.Range("E25").Value = "'  -2.74%  "
$ws.Range("D26").Value = "'0.179"
$ws.Range("E26").Value = "'  -1.42%  "
$ws.Range("D27").Value = "'27.00"
$ws.Range("E27").Value = "'  +0.35%  "
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "'  +0.16%  "
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").Value = "'7.47"
$ws.Range("E29").Value = "'  +2.46%  "
$ws.Range("E30").Value = "'  +1.76%  "
$ws.Range("D31").Value = "'6.54"
$ws.Range("E31").Value = "'  +8.92%  "
$ws.Range("E32").Value = "'  -2.93%  "
$ws.Range("D33").Value = "'36.00"
$ws.Range("E33").Value = "'  -6.58%  "
$ws.Range("D34").Value = "'2.18"
$ws.Range("E34").Value = "'  +15.83%  "
$ws.Range("D35").Value = "'51.18"
$ws.Range("E35").Value = "'  -1.65%  "
$ws.Range("E36").Value = "'  -0.16%  "
$ws.Range("E37").Value = "'  +0.05%  "
$ws.Range("D38").Value = "'3.32"
$ws.Range("E38").Value = "'  +3.60%  "
$ws.Range("D39").Value = "'2.82"
$ws.Range("E39").Value = "'  +3.54%  "
$ws.Range("E40").Value = "'  -1.99%  "
$ws.Range("D41").Value = "'17.51"
$ws.Range("E41").Value = "'  -3.89%  "
$ws.Range("E42").Value = "'  -1.83%  "
$ws.Range("D43").Value = "'23.41"
$ws.Range("E43").Value = "'  +2.49%  "
$ws.Range("D44").Value = "'125.28"
$ws.Range("E44").Value = "'  +4.15%  "
$ws.Range("D45").Value = "'2.18"
$ws.Range("E45").Value = "'  +0.86%  "
$ws.Range("D46").Value = "'2.127.53"
$ws.Range("E46").Value = "'  +0.26%  "
$ws.Range("E47").Value = "'  -2.63%  "
$ws.Range("D48").Value = "'2.34"
$ws.Range("E48").Value = "'  -6.78%  "
$ws.Range("D49").Value = "'0.246"
$ws.Range("E49").Value = "'  +0.01%  "
$ws.Range("E50").Value = "'  +3.07%  "
$ws.Range("D51").Value = "'0.906"
$ws.Range("E51").Value = "'  +1.29%  "
